# "affichage graphique dans Excel"
# - Move the two stray section-header labels from column C to column B
#   (row 16 "Contenu du stage", row 25 "Type entreprise") so they line up
#   with the other section headers (B6, B8, B10) used as chart series names.
# - Insert three pie charts (Lieu du stage / Contenu du stage / Type du
#   stage) driven by the D/E columns, stacked down column I, and wire the
#   worksheet drawing relationship.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Re-home the section-header labels: C16 -> B16, C25 -> B25
# ---------------------------------------------------------------------
$c16 = $ws.Range("C16").Value2
$ws.Range("B16").Value2 = $c16
$ws.Range("C16").ClearContents()

$c25 = $ws.Range("C25").Value2
$ws.Range("B25").Value2 = $c25
$ws.Range("C25").ClearContents()

# ---------------------------------------------------------------------
# Helper geometry: charts are stacked in the I2:Q.. block, one above the
# other, each spanning 15 rows (same footprint the workbook author used).
# ---------------------------------------------------------------------
$left = $ws.Cells.Item(2, 9).Left
$width = $ws.Cells.Item(2, 17).Left - $left

$orange = 39423  # RGB(255,153,0) -> FF9900

# ---------------------------------------------------------------------
# 2) Chart 1 - "Lieu du stage" (B10 / D10:D14 / E10:E14)
# ---------------------------------------------------------------------
$top1 = $ws.Cells.Item(2, 9).Top
$height1 = $ws.Cells.Item(16, 9).Top - $top1
$co1 = $ws.ChartObjects().Add($left, $top1, $width, $height1)
$co1.Name = "Chart 1"
$chart1 = $co1.Chart
$chart1.ChartType = 5
$ser1 = $chart1.SeriesCollection().NewSeries()
$ser1.Name = "=Worksheet!`$B`$10"
$ser1.XValues = $ws.Range("D10:D14")
$ser1.Values = $ws.Range("E10:E14")
$chart1.HasTitle = $true
$chart1.ChartTitle.Text = "Lieu du stage"
$chart1.HasLegend = $true
$chart1.Legend.Position = -4152
$ser1.HasDataLabels = $true
$dl1 = $ser1.DataLabels()
$dl1.ShowValue = $true
$dl1.ShowPercentage = $true
$dl1.ShowCategoryName = $false
$dl1.ShowSeriesName = $false
$dl1.ShowLegendKey = $false
$pt1 = $ser1.Points(4)
$pt1.Interior.Color = $orange
$pt1.Format.Fill.ForeColor.RGB = $orange

# ---------------------------------------------------------------------
# 3) Chart 2 - "Contenu du stage" (B16 / D16:D23 / E16:E23)
# ---------------------------------------------------------------------
$top2 = $ws.Cells.Item(18, 9).Top
$height2 = $ws.Cells.Item(32, 9).Top - $top2
$co2 = $ws.ChartObjects().Add($left, $top2, $width, $height2)
$co2.Name = "Chart 2"
$chart2 = $co2.Chart
$chart2.ChartType = 5
$ser2 = $chart2.SeriesCollection().NewSeries()
$ser2.Name = "=Worksheet!`$B`$16"
$ser2.XValues = $ws.Range("D16:D23")
$ser2.Values = $ws.Range("E16:E23")
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = "Contenu du stage"
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152
$ser2.HasDataLabels = $true
$dl2 = $ser2.DataLabels()
$dl2.ShowValue = $true
$dl2.ShowPercentage = $true
$dl2.ShowCategoryName = $false
$dl2.ShowSeriesName = $false
$dl2.ShowLegendKey = $false
$pt2 = $ser2.Points(4)
$pt2.Interior.Color = $orange
$pt2.Format.Fill.ForeColor.RGB = $orange

# ---------------------------------------------------------------------
# 4) Chart 3 - title "Type du stage", series still sourced from B25
#    ("Type entreprise") / D25:D28 / E25:E28
# ---------------------------------------------------------------------
$top3 = $ws.Cells.Item(34, 9).Top
$height3 = $ws.Cells.Item(50, 9).Top - $top3
$co3 = $ws.ChartObjects().Add($left, $top3, $width, $height3)
$co3.Name = "Chart 3"
$chart3 = $co3.Chart
$chart3.ChartType = 5
$ser3 = $chart3.SeriesCollection().NewSeries()
$ser3.Name = "=Worksheet!`$B`$25"
$ser3.XValues = $ws.Range("D25:D28")
$ser3.Values = $ws.Range("E25:E28")
$chart3.HasTitle = $true
$chart3.ChartTitle.Text = "Type du stage"
$chart3.HasLegend = $true
$chart3.Legend.Position = -4152
$ser3.HasDataLabels = $true
$dl3 = $ser3.DataLabels()
$dl3.ShowValue = $true
$dl3.ShowPercentage = $true
$dl3.ShowCategoryName = $false
$dl3.ShowSeriesName = $false
$dl3.ShowLegendKey = $false
$pt3 = $ser3.Points(4)
$pt3.Interior.Color = $orange
$pt3.Format.Fill.ForeColor.RGB = $orange
